$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 505, shifting existing rows 505:542 down to 506:543.
$ws.Rows.Item(505).Insert()

# Populate the new row 505 with the new weekly record.
$ws.Range("A505").Value = 3
$ws.Range("B505").Value = "Femacal de La Calera"
$ws.Range("C505").Value = "Coquimbo"
$ws.Range("D505").Value = 44826
$ws.Range("E505").Value = 5
$ws.Range("F505").Value = 100112021
$ws.Range("G505").Value = "Ají"
$ws.Range("H505").Value = "Inferno"
$ws.Range("I505").Value = "Primera"
$ws.Range("J505").Value = 76
$ws.Range("K505").Value = 16000
$ws.Range("L505").Value = 17000
$ws.Range("M505").Value = 16500
$ws.Range("N505").Value = "$/caja 10 kilos"
$ws.Range("O505").Value = "Región de Arica y Parinacota"
$ws.Range("P505").Value = 1650
$ws.Range("Q505").Value = 10
$ws.Range("R505").Value = "Hortaliza"
